# BAARD uses fixed training set - update evaluation metrics in the
# "breastcancer" worksheet to reflect the re-run with a fixed training set.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("breastcancer")

# Row 4 (apgd, 0.05)
$ws.Range("C4").Value = 89.24731182795698
$ws.Range("D4").Value = 90.32258064516128

# Row 5 (apgd, 0.2)
$ws.Range("C5").Value = 90.32258064516128
$ws.Range("D5").Value = 90.32258064516128
$ws.Range("H5").Value = 13.97849462365591

# Row 6 (apgd, 0.4)
$ws.Range("H6").Value = 13.97849462365591

# Row 7 (apgd, 1)
$ws.Range("G7").Value = 11.82795698924731

# Row 8 (apgd2, 0.4)
$ws.Range("C8").Value = 74.19354838709677
$ws.Range("D8").Value = 79.56989247311827
$ws.Range("H8").Value = 13.97849462365591

# Row 9 (apgd2, 1)
$ws.Range("C9").Value = 60.21505376344086
$ws.Range("D9").Value = 60.21505376344086
$ws.Range("H9").Value = 13.97849462365591

# Row 10 (apgd2, 2)
$ws.Range("C10").Value = 90.32258064516128
$ws.Range("D10").Value = 90.32258064516128
$ws.Range("G10").Value = 11.82795698924731
$ws.Range("H10").Value = 13.97849462365591

# Row 11 (apgd2, 3)
$ws.Range("G11").Value = 11.82795698924731
$ws.Range("H11").Value = 13.97849462365591

# Row 12 (boundary, 0.3)
$ws.Range("C12").Value = 87.09677419354837
$ws.Range("G12").Value = 11.82795698924731
$ws.Range("H12").Value = 13.97849462365591

# Row 13 (cw2, 0)
$ws.Range("C13").Value = 72.04301075268819
$ws.Range("D13").Value = 94.6236559139785
$ws.Range("G13").Value = 11.82795698924731

# Row 14 (cw2, 5)
$ws.Range("C14").Value = 52.68817204301075
$ws.Range("D14").Value = 52.68817204301075
$ws.Range("G14").Value = 11.82795698924731
$ws.Range("H14").Value = 13.97849462365591

# Row 15 (cw2, 10)
$ws.Range("C15").Value = 86.02150537634407
$ws.Range("D15").Value = 86.02150537634407
$ws.Range("G15").Value = 11.82795698924731

# Row 16 (deepfool, 1e-06)
$ws.Range("C16").Value = 84.94623655913979

# Row 17 (fgsm, 0.05)
$ws.Range("C17").Value = 89.24731182795698
$ws.Range("D17").Value = 89.24731182795698
$ws.Range("G17").Value = 11.82795698924731
$ws.Range("H17").Value = 13.97849462365591

# Row 18 (fgsm, 0.2)
$ws.Range("C18").Value = 93.54838709677419
$ws.Range("G18").Value = 11.82795698924731
$ws.Range("H18").Value = 13.97849462365591

# Row 19 (fgsm, 0.4)
$ws.Range("H19").Value = 13.97849462365591

# Row 20 (fgsm, 1)
$ws.Range("H20").Value = 13.97849462365591
